$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("B2").Value = 509.291461187893
$ws.Range("C2").Value = 355.3664490273144
$ws.Range("D2").Value = 304.5855631686791
$ws.Range("E2").Value = 280.1723854596627
$ws.Range("B3").Value = 598.4954603342391
$ws.Range("C3").Value = 418.9857197695599
$ws.Range("D3").Value = 357.3558453484139
$ws.Range("E3").Value = 331.565737373865
$ws.Range("B4").Value = 572.6393071452801
$ws.Range("C4").Value = 402.1265504482542
$ws.Range("D4").Value = 343.9891792172732
$ws.Range("E4").Value = 319.5669083563943
$ws.Range("B5").Value = 393.4905759255764
$ws.Range("C5").Value = 276.213026956218
$ws.Range("D5").Value = 232.832744098615
$ws.Range("E5").Value = 218.6987064554309
$ws.Range("B6").Value = 350.5477972904769
$ws.Range("C6").Value = 244.3955655368507
$ws.Range("D6").Value = 207.1470503427233
$ws.Range("E6").Value = 191.9618107489409
$ws.Range("B7").Value = 36.65128319507506
$ws.Range("C7").Value = 25.57506905940752
$ws.Range("D7").Value = 21.91258491892667
$ws.Range("E7").Value = 20.24933835187164
$ws.Range("B8").Value = 1996.711446961557
$ws.Range("C8").Value = 1393.477163308715
$ws.Range("D8").Value = 1204.286864990044
$ws.Range("E8").Value = 1103.377142008882
$ws.Range("B9").Value = 502.0824162856143
$ws.Range("C9").Value = 351.7372728805154
$ws.Range("D9").Value = 299.8303694770588
$ws.Range("E9").Value = 278.4862426358364
$ws.Range("B10").Value = 204.8215920687625
$ws.Range("C10").Value = 143.3785024642036
$ws.Range("D10").Value = 125.7069016271475
$ws.Range("E10").Value = 115.1558680344456
$ws.Range("B11").Value = 37.72954408524613
$ws.Range("C11").Value = 25.07765809233451
$ws.Range("D11").Value = 21.71277507041724
$ws.Range("E11").Value = 21.44469988760394
$ws.Range("B12").Value = 87.10810053658919
$ws.Range("C12").Value = 64.69472673945423
$ws.Range("D12").Value = 56.97064824144755
$ws.Range("E12").Value = 51.6504788309662
$ws.Range("B13").Value = 114.7584831027659
$ws.Range("C13").Value = 79.67824661358124
$ws.Range("D13").Value = 70.60540822888949
$ws.Range("E13").Value = 65.46501560483037
